# "Update database" - Update Tracking Process.xlsx
#
# On the "Database" sheet, mark the "Review Database" (row 5) and
# "Generate Script from Draft" (row 6) tasks as 100% completed by
# setting their "% Completed" (column D) value to 1, matching the
# value already used for the first (100%-complete) task in row 4.
$wb = $excel.ActiveWorkbook

$wsDatabase = $wb.Worksheets.Item("Database")
$wsDatabase.Range("D5").Value = 1
$wsDatabase.Range("D6").Value = 1

# Reflect the user's navigation: they were last looking at the
# "Document" sheet with cell B27 selected ...
$wsDocument = $wb.Worksheets.Item("Document")
[void]$wsDocument.Activate()
[void]$wsDocument.Range("B27").Select()

# ... then moved to the "Database" sheet (now the active / selected tab)
# with cell G7 selected, which is where the workbook was left.
[void]$wsDatabase.Activate()
[void]$wsDatabase.Range("G7").Select()
